$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 3-7 as part of repulled data / mean calculation
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -3
